# Refresh the cryptos list values to match the latest GitHub Actions scrape.
# Column D ("Price") holds text like "1.00" / "96.203.15" / "3.922.62" that Excel
# would otherwise auto-coerce into a Double (dropping the text formatting /
# multi-dot thousands separators), so those assignments are apostrophe-prefixed
# to force a literal text entry, matching the source workbook (t="inlineStr").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'" + '96.203.15'
$ws.Range("E2").Value = '  -2.33%  '

$ws.Range("D3").Value = "'" + '3.308.87'
$ws.Range("E3").Value = '  -4.52%  '

$ws.Range("E4").Value = '  +0.10%  '

$ws.Range("D5").Value = "'" + '245.53'
$ws.Range("E5").Value = '  -5.55%  '

$ws.Range("D6").Value = "'" + '646.62'
$ws.Range("E6").Value = '  -3.52%  '

$ws.Range("D7").Value = "'" + '1.32'
$ws.Range("E7").Value = '  -14.68%  '

$ws.Range("D8").Value = "'" + '0.405'

$ws.Range("D9").Value = "'" + '1.00'
$ws.Range("E9").Value = '  +0.15%  '

$ws.Range("D10").Value = "'" + '0.957'
$ws.Range("E10").Value = '  -13.19%  '

$ws.Range("D11").Value = "'" + '3.308.61'
$ws.Range("E11").Value = '  -4.45%  '

$ws.Range("E12").Value = '  -7.27%  '

$ws.Range("D13").Value = "'" + '39.12'
$ws.Range("E13").Value = '  -9.36%  '

$ws.Range("D14").Value = "'" + '95.977.06'
$ws.Range("E14").Value = '  -1.96%  '

$ws.Range("E15").Value = '  -5.61%  '

$ws.Range("B16").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C16").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D16").Value = "'" + '3.922.62'
$ws.Range("E16").Value = '  -4.64%  '

$ws.Range("B17").Value = 'ShibaInu'
$ws.Range("C17").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D17").Value = "'" + '0.0000246'
$ws.Range("E17").Value = '  -9.00%  '

$ws.Range("D18").Value = "'" + '8.37'
$ws.Range("E18").Value = '  -5.04%  '

$ws.Range("D19").Value = "'" + '3.295.30'
$ws.Range("E19").Value = '  -4.80%  '

$ws.Range("D20").Value = "'" + '16.50'
$ws.Range("E20").Value = '  -7.73%  '

$ws.Range("D21").Value = "'" + '10.26'
$ws.Range("E21").Value = '  -8.80%  '

$ws.Range("D22").Value = "'" + '489.37'
$ws.Range("E22").Value = '  -6.17%  '

$ws.Range("E23").Value = '  -7.98%  '

$ws.Range("D24").Value = "'" + '0.460'
$ws.Range("E24").Value = '  -10.17%  '

$ws.Range("D25").Value = "'" + '0.0000194'
$ws.Range("E25").Value = '  -10.43%  '

$ws.Range("D26").Value = "'" + '6.32'
$ws.Range("E26").Value = '  -2.16%  '

$ws.Range("D27").Value = "'" + '91.21'
$ws.Range("E27").Value = '  -11.49%  '

$ws.Range("D28").Value = "'" + '11.76'
$ws.Range("E28").Value = '  -9.74%  '

$ws.Range("D29").Value = "'" + '3.483.72'
$ws.Range("E29").Value = '  -4.49%  '

$ws.Range("E30").Value = '  +0.13%  '

$ws.Range("E31").Value = '  -14.05%  '

$ws.Range("D32").Value = "'" + '10.54'
$ws.Range("E32").Value = '  -10.90%  '

$ws.Range("E33").Value = '  -7.55%  '

$ws.Range("D34").Value = "'" + '2.40'
$ws.Range("E34").Value = '  +6.48%  '

$ws.Range("D35").Value = "'" + '1.00'
$ws.Range("E35").Value = '  -0.49%  '

$ws.Range("D36").Value = "'" + '0.531'
$ws.Range("E36").Value = '  -11.10%  '

$ws.Range("D37").Value = "'" + '27.45'
$ws.Range("E37").Value = '  -9.89%  '

$ws.Range("D38").Value = "'" + '1.43'
$ws.Range("E38").Value = '  +0.18%  '

$ws.Range("D39").Value = "'" + '7.37'
$ws.Range("E39").Value = '  -8.71%  '

$ws.Range("E40").Value = '  -0.06%  '

$ws.Range("D41").Value = "'" + '0.147'
$ws.Range("E41").Value = '  -8.27%  '

$ws.Range("D42").Value = "'" + '494.80'
$ws.Range("E42").Value = '  -7.33%  '

$ws.Range("D43").Value = "'" + '24.45'
$ws.Range("E43").Value = '  -1.23%  '

$ws.Range("E44").Value = '  -2.25%  '

$ws.Range("D45").Value = "'" + '0.806'
$ws.Range("E45").Value = '  -7.08%  '

$ws.Range("D46").Value = "'" + '0.0395'
$ws.Range("E46").Value = '  -12.36%  '

$ws.Range("D47").Value = "'" + '8.18'
$ws.Range("E47").Value = '  -4.56%  '

$ws.Range("E48").Value = '  -1.66%  '

$ws.Range("D49").Value = "'" + '1.57'
$ws.Range("E49").Value = '  -2.67%  '

$ws.Range("D50").Value = "'" + '52.15'
$ws.Range("E50").Value = '  +0.51%  '

$ws.Range("D51").Value = "'" + '3.09'
$ws.Range("E51").Value = '  -10.11%  '
